$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.781.73"
$ws.Range("D3").Value = "1.655.59"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3804"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3630"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.19"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.257"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08231"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.74"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.546"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.473"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001242"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "1.656.41"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.80"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06980"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.804"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.80"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("D24").Value = "23.786.66"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.561"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.083"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.35"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.38"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.232"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.57"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "1.840.86"
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.920"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.193"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.083"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.85"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02835"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2530"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.137"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08821"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07131"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +11.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7089"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.345"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.03"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6575"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.338"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.966"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07964"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.18"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("E51").Value = "  +0.74%  "
